# Personnels logic JUST DRAFT
# Adds a new column L with a literal 1 marker value on the two existing
# data rows (row 1 = header-ish first record, row 2 = second record),
# matching the author's manual addition of a new field/flag column next
# to the existing "К" (marital status) column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("L1").Value = 1
$ws.Range("L2").Value = 1

# Mirror the author's final on-screen selection as closely as the
# object model allows: they had clicked L2 and then ctrl-selected the
# whole of column N, leaving N1 as the active cell.
$ws.Columns.Item(14).Select()
